# System Setup4: sydefault, Inv comcod
#
# Update the "Edit_ScrapReasonCode" sheet:
#  - B2 text: "Scrap Reason Code 1" -> "Scrap Reason Code Update"
#  - Column B gets an explicit best-fit width so the updated text is fully visible
#  - Active selection moves to B5

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Edit_ScrapReasonCode")
$ws.Activate()

$ws.Range("B2").Value = "Scrap Reason Code Update"

$ws.Columns.Item(2).ColumnWidth = 17.3

$ws.Range("B5").Select()
